$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text) ---
$src = $ws.Range("A8")
$full = $src.Text
$idx = $full.IndexOf("48")
$chars = $src.Characters($idx + 1, 2)
$chars.Text = "49"

$src2 = $ws.Range("C9")
$full2 = $src2.Text
$idx2 = $full2.IndexOf("11/28/2022")
$chars2 = $src2.Characters($idx2 + 1, 10)
$chars2.Text = "12/5/2022"
$full2b = $src2.Text
$idx2b = $full2b.IndexOf("12/4/2022")
$chars2b = $src2.Characters($idx2b + 1, 9)
$chars2b.Text = "12/11/2022"

# --- Data grid updates (rows 14-29) ---
# Row 14
$ws.Range("N14").Value = -90.625

# Row 15
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("N23").Copy($ws.Range("E15"))
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -66.666666666666
$ws.Range("L15").Value = 36.363636363636

# Row 16
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 44
$ws.Range("G16").Value = 56
$ws.Range("H16").Value = -21.428571428571
$ws.Range("I16").Value = 700
$ws.Range("J16").Value = 493
$ws.Range("K16").Value = 41.987829614604
$ws.Range("L16").Value = 87.667560321715
$ws.Range("M16").Value = 41.700404858299
$ws.Range("N16").Value = -64.084145715751

# Row 17
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = 37.5
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = 16.363636363636
$ws.Range("I17").Value = 761
$ws.Range("J17").Value = 672
$ws.Range("K17").Value = 13.244047619047
$ws.Range("L17").Value = 19.654088050314
$ws.Range("M17").Value = 54.989816700611
$ws.Range("N17").Value = -15.350389321468

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 14
$ws.Range("I18").Value = 289
$ws.Range("J18").Value = 232
$ws.Range("K18").Value = 24.568965517241
$ws.Range("L18").Value = 31.363636363636
$ws.Range("M18").Value = -22.10242587601
$ws.Range("N18").Value = -84.25068119891

# Row 19
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = 13.636363636363
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 92
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 992
$ws.Range("J19").Value = 758
$ws.Range("K19").Value = 30.870712401055
$ws.Range("L19").Value = 55
$ws.Range("M19").Value = 83.025830258302
$ws.Range("N19").Value = 33.333333333333

# Row 20
$ws.Range("C20").Value = 21
$ws.Range("E20").Value = 90.90909090909
$ws.Range("F20").Value = 63
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = 46.511627906976
$ws.Range("I20").Value = 521
$ws.Range("J20").Value = 420
$ws.Range("K20").Value = 24.047619047619
$ws.Range("L20").Value = 136.818181818182
$ws.Range("M20").Value = 125.541125541126
$ws.Range("N20").Value = -73.199588477366

# Row 21
$ws.Range("C21").Value = 84
$ws.Range("D21").Value = 68
$ws.Range("E21").Value = 23.529411764705
$ws.Range("F21").Value = 279
$ws.Range("G21").Value = 267
$ws.Range("H21").Value = 4.494382022471
$ws.Range("I21").Value = 3314
$ws.Range("J21").Value = 2633
$ws.Range("K21").Value = 25.864033421952
$ws.Range("L21").Value = 55.513843266072
$ws.Range("M21").Value = 52.297794117647
$ws.Range("N21").Value = -55.83100093296

# Row 22
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 5.263157894736

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 34.782608695652
$ws.Range("I23").Value = 320
$ws.Range("J23").Value = 254
$ws.Range("K23").Value = 25.984251968503
$ws.Range("L23").Value = 15.107913669064
$ws.Range("M23").Value = 47.465437788018

# Row 24
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 159
$ws.Range("G24").Value = 166
$ws.Range("H24").Value = -4.216867469879
$ws.Range("I24").Value = 1968
$ws.Range("J24").Value = 1385
$ws.Range("K24").Value = 42.093862815884
$ws.Range("L24").Value = 48.640483383685
$ws.Range("M24").Value = 35.817805383022

# Row 25
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = -53.333333333333
$ws.Range("F25").Value = 72
$ws.Range("G25").Value = 87
$ws.Range("H25").Value = -17.241379310344
$ws.Range("I25").Value = 1063
$ws.Range("J25").Value = 951
$ws.Range("K25").Value = 11.777076761303
$ws.Range("L25").Value = 14.424111948331
$ws.Range("M25").Value = -27.834351663272

# Row 26
$ws.Range("C15").Copy($ws.Range("D26"))
$ws.Range("N23").Copy($ws.Range("E26"))
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = -83.333333333333
$ws.Range("L26").Value = 32

# Row 27
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 68
$ws.Range("J27").Value = 97
$ws.Range("K27").Value = -29.896907216494
$ws.Range("L27").Value = 7.936507936507

# Row 28
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("N23").Copy($ws.Range("E28"))
$ws.Range("C15").Copy($ws.Range("F28"))
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -72.48322147651

# Row 29
$ws.Range("C15").Copy($ws.Range("D29"))
$ws.Range("N23").Copy($ws.Range("E29"))
$ws.Range("C15").Copy($ws.Range("F29"))
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -74.814814814814
